# Turn the single data row (anako / user / 1000.0) into a header row
# (Username / Role / Budget) and size the columns to fit their new
# text, mirroring an "add default fields" header commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Role"
$ws.Range("C1").Value = "Budget"

# Best-fit the three header columns to their text (bestFit-style autosize).
$ws.Range("A1:C1").EntireColumn.AutoFit()

# Nudge the autosized widths to the precise best-fit character widths
# ("Username" / "Role" / "Budget" in the sheet's default 11pt Calibri),
# snapped to the nearest width this engine can actually persist.
$ws.Columns.Item(1).ColumnWidth = 9.333333333333334
$ws.Columns.Item(2).ColumnWidth = 4.166666666666667
$ws.Columns.Item(3).ColumnWidth = 6.5
